# dias_franjas.xlsx edit:
#   - "franjas" renamed to "franjas v1"
#   - new sheet "franjas v2" inserted right after "franjas v1" (a time-grid /
#     teacher-availability schedule sheet)
#   - selection markers updated on "franjas v1" and the new "franjas v2"
#
# disp1 / disp2 are left untouched (they just shift from rId2/rId3 to
# rId3/rId4, which Excel/the engine renumbers automatically when the new
# sheet is inserted).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet and insert the new one right after it.
# ---------------------------------------------------------------------------
$v1 = $wb.Worksheets.Item("franjas")
$v1.Name = "franjas v1"

$v2 = $wb.Worksheets.Add($null, $v1)
$v2.Name = "franjas v2"

# ---------------------------------------------------------------------------
# 2. Row 1 / Row 2 : the 15-minute time-slot header grid, columns B..BI.
#    (There's a genuine jump from 2100 to 2215 near the end of the source
#    data, so these are listed literally rather than generated by a +15
#    loop.)
# ---------------------------------------------------------------------------
$row1 = @(700,715,730,745,800,815,830,845,900,915,930,945,1000,1015,1030,1045,1100,1115,1130,1145,1200,1215,1230,1245,1300,1315,1330,1345,1400,1415,1430,1445,1500,1515,1530,1545,1600,1615,1630,1645,1700,1715,1730,1745,1800,1815,1830,1845,1900,1915,1930,1945,2000,2015,2030,2045,2100,2215,2230,2245)
$row2 = @(715,730,745,800,815,830,845,900,915,930,945,1000,1015,1030,1045,1100,1115,1130,1145,1200,1215,1230,1245,1300,1315,1330,1345,1400,1415,1430,1445,1500,1515,1530,1545,1600,1615,1630,1645,1700,1715,1730,1745,1800,1815,1830,1845,1900,1915,1930,1945,2000,2015,2030,2045,2100,2215,2230,2245,2300)

$headerGrid = New-Object 'object[,]' 2,60
for ($i = 0; $i -lt 60; $i++) {
    $headerGrid[0,$i] = $row1[$i]
    $headerGrid[1,$i] = $row2[$i]
}
$v2.Range("B1:BI2").Value = $headerGrid
$v2.Range("B1:BI2").Style = "Normal 2"

# ---------------------------------------------------------------------------
# 3. Rows 3-8 : single-letter weekday labels in column A (M,T,W,R,F,S -
#    these reuse shared-string entries already present in the workbook).
# ---------------------------------------------------------------------------
$dayLabels = @("M","T","W","R","F","S")
for ($i = 0; $i -lt $dayLabels.Length; $i++) {
    $v2.Cells.Item(3 + $i, 1).Value = $dayLabels[$i]
}
$v2.Range("A3:A8").Style = "Normal 2"

# ---------------------------------------------------------------------------
# 4. Rows 10-69 : contiguous 15-minute (hi, hf) pairs in columns D, E,
#    from 700-715 up to 2145-2200.
# ---------------------------------------------------------------------------
function Get-NextQuarterHour($t) {
    $h = [math]::Floor($t / 100)
    $m = $t % 100
    $m += 15
    if ($m -ge 60) {
        $m -= 60
        $h += 1
    }
    return ($h * 100) + $m
}

$deGrid = New-Object 'object[,]' 60,2
$t = 700
for ($i = 0; $i -lt 60; $i++) {
    $next = Get-NextQuarterHour $t
    $deGrid[$i,0] = $t
    $deGrid[$i,1] = $next
    $t = $next
}
$v2.Range("D10:E69").Value = $deGrid
$v2.Range("D10:E69").Style = "Normal 2"

# ---------------------------------------------------------------------------
# 5. Column widths for the new sheet (best-effort match of the authored
#    bestFit widths; the COM layer only keeps 2 decimal "character" digits).
# ---------------------------------------------------------------------------
$v2.Columns.Item(1).ColumnWidth = 2
$v2.Range("B1:C1").EntireColumn.ColumnWidth = 3.33
$v2.Range("D1:E1").EntireColumn.ColumnWidth = 4.33
$v2.Range("F1:L1").EntireColumn.ColumnWidth = 3.33
$v2.Range("M1:BI1").EntireColumn.ColumnWidth = 4.33

# ---------------------------------------------------------------------------
# 6. Selections: "franjas v1" keeps the grid selected at K27, the new
#    "franjas v2" (now the active tab) is selected at C9.
# ---------------------------------------------------------------------------
$v1.Range("K27").Select()
$v2.Range("C9").Select()
